$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Nid2"
$ws.Range("C2").Value = "Col13a1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 16.633059
$ws.Range("H2").Value = 49.899177
$ws.Range("I2").Value = 0.2141889211266894
$ws.Range("J2").Value = 0.2141889211266894
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.3326186666666667
$ws.Range("N2").Value = 0.9978560000000001
$ws.Range("O2").Value = 0.3213922220876632
$ws.Range("P2").Value = 0.3213922220876632
$ws.Range("Q2").Value = 5.532465907168
$ws.Range("R2").Value = 49.79219316451201
$ws.Range("S2").Value = 0.06883865330746594
$ws.Range("T2").Value = 0.06883865330746594

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Nid2"
$ws.Range("C3").Value = "Col13a1"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 16.633059
$ws.Range("H3").Value = 49.899177
$ws.Range("I3").Value = 0.2141889211266894
$ws.Range("J3").Value = 0.2141889211266894
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.702312
$ws.Range("N3").Value = 2.106936
$ws.Range("O3").Value = 0.6786077779123368
$ws.Range("P3").Value = 0.6786077779123368
$ws.Range("Q3").Value = 11.681596932408
$ws.Range("R3").Value = 105.134372391672
$ws.Range("S3").Value = 0.1453502678192235
$ws.Range("T3").Value = 0.1453502678192235

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Nid2"
$ws.Range("C4").Value = "Col13a1"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 50.79415266666667
$ws.Range("H4").Value = 152.382458
$ws.Range("I4").Value = 0.6540916351717195
$ws.Range("J4").Value = 0.6540916351717195
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.3326186666666667
$ws.Range("N4").Value = 0.9978560000000001
$ws.Range("O4").Value = 0.3213922220876632
$ws.Range("P4").Value = 0.3213922220876632
$ws.Range("Q4").Value = 16.89508333444978
$ws.Range("R4").Value = 152.055750010048
$ws.Range("S4").Value = 0.2102199640767921
$ws.Range("T4").Value = 0.2102199640767921

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Nid2"
$ws.Range("C5").Value = "Col13a1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 50.79415266666667
$ws.Range("H5").Value = 152.382458
$ws.Range("I5").Value = 0.6540916351717195
$ws.Range("J5").Value = 0.6540916351717195
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.702312
$ws.Range("N5").Value = 2.106936
$ws.Range("O5").Value = 0.6786077779123368
$ws.Range("P5").Value = 0.6786077779123368
$ws.Range("Q5").Value = 35.673342947632
$ws.Range("R5").Value = 321.060086528688
$ws.Range("S5").Value = 0.4438716710949275
$ws.Range("T5").Value = 0.4438716710949275

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Nid2"
$ws.Range("C6").Value = "Col13a1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 10.22880766666667
$ws.Range("H6").Value = 30.686423
$ws.Range("I6").Value = 0.1317194437015911
$ws.Range("J6").Value = 0.1317194437015911
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.3326186666666667
$ws.Range("N6").Value = 0.9978560000000001
$ws.Range("O6").Value = 0.3213922220876632
$ws.Range("P6").Value = 0.3213922220876632
$ws.Range("Q6").Value = 3.402292367676444
$ws.Range("R6").Value = 30.620631309088
$ws.Range("S6").Value = 0.04233360470340521
$ws.Range("T6").Value = 0.04233360470340521

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Nid2"
$ws.Range("C7").Value = "Col13a1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 10.22880766666667
$ws.Range("H7").Value = 30.686423
$ws.Range("I7").Value = 0.1317194437015911
$ws.Range("J7").Value = 0.1317194437015911
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.702312
$ws.Range("N7").Value = 2.106936
$ws.Range("O7").Value = 0.6786077779123368
$ws.Range("P7").Value = 0.6786077779123368
$ws.Range("Q7").Value = 7.183814369992001
$ws.Range("R7").Value = 64.654329329928
$ws.Range("S7").Value = 0.08938583899818586
$ws.Range("T7").Value = 0.08938583899818586
